# Insert a new data record at row 51 of the "Macroferia Regional de Talca -
# Cilantro" sheet. All existing records from row 51 down to row 121 shift
# down by one row (to rows 52..122), and the sheet's used range grows from
# A1:R121 to A1:R122.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 51..121 down to 52..122, leaving a blank row 51 to fill in.
$ws.Rows(51).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(51, 1).Value  = 5
$ws.Cells.Item(51, 2).Value  = 'Macroferia Regional de Talca'
$ws.Cells.Item(51, 3).Value  = 'Maule'
$ws.Cells.Item(51, 4).Value  = 45174
$ws.Cells.Item(51, 5).Value  = 7
$ws.Cells.Item(51, 6).Value  = 100112040
$ws.Cells.Item(51, 7).Value  = 'Cilantro'
$ws.Cells.Item(51, 8).Value  = 'Sin especificar'
$ws.Cells.Item(51, 9).Value  = 'Primera'
$ws.Cells.Item(51, 10).Value = 200
$ws.Cells.Item(51, 11).Value = 8000
$ws.Cells.Item(51, 12).Value = 8000
$ws.Cells.Item(51, 13).Value = 8000
$ws.Cells.Item(51, 14).Value = '$/caja 36 atados'
$ws.Cells.Item(51, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(51, 16).Value = 222
$ws.Cells.Item(51, 17).Value = 36
$ws.Cells.Item(51, 18).Value = 'Hortaliza'
